# Insert a new weekly price record as row 327 on the sheet, pushing the
# existing rows 327:339 down to 328:340 (same as Excel's "Insert Sheet Rows").
# The new row carries the same Mercado/Region/Categoria context as the rest
# of the block, with its own date, price and origin figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(327).Insert()

$ws.Cells.Item(327, 1).Value = 10
$ws.Cells.Item(327, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(327, 3).Value = "La Araucanía"
$ws.Cells.Item(327, 4).Value = 44509
$ws.Cells.Item(327, 5).Value = 9
$ws.Cells.Item(327, 6).Value = 100112043
$ws.Cells.Item(327, 7).Value = "Pepino ensalada"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 200
$ws.Cells.Item(327, 11).Value = 10000
$ws.Cells.Item(327, 12).Value = 10000
$ws.Cells.Item(327, 13).Value = 10000
$ws.Cells.Item(327, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(327, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(327, 16).Value = 167
$ws.Cells.Item(327, 17).Value = 60
$ws.Cells.Item(327, 18).Value = "Hortaliza"
